# Combine stations 5 + 5b for more figures
#
# The "Littorina TBT" sheet gets a new "Station_name" column inserted right
# after "TBT_flag" (i.e. becomes the new column D), with the descriptive
# station name for each of the five station rows. All of the existing
# year-columns (2005..2018) shift one column to the right to make room.
#
# The other three sheets (Littorina ISI / sterile / PRL) are untouched by
# this change - their cell content stays exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Littorina TBT")

# Insert a new, blank column at D - this shifts the old D:L (2005..2018)
# data over to E:M, matching the diff.
$ws.Columns.Item(4).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 4).Value = "Station_name"

# Station descriptions for each data row (station numbers are in column B).
$ws.Cells.Item(2, 4).Value = "Reference station (5.5 km)"
$ws.Cells.Item(3, 4).Value = "Outer Vikkilen (2.5 km)"
$ws.Cells.Item(4, 4).Value = "Skjeviga (0.1 km)"
$ws.Cells.Item(5, 4).Value = "Shipyard (0 km)"
$ws.Cells.Item(6, 4).Value = "Inner Vikkilen (0.5 km)"
